$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format so numeric-looking strings (prices, percentages)
# are stored verbatim instead of being parsed/rounded as numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.375.91'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.59%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.686.69'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.08%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '687.79'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.69'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -5.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.685.42'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.09%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -5.45%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -8.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.38'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.20%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -7.15%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.63%  '
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.306.74'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.17%  '
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '33.05'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -7.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.684.70'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.382.98'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.113'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.47%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -7.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.53'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -8.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '476.43'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -7.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.96'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.30%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -7.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.85'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.831.11'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000129'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -9.29%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.23'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -7.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.37'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -9.12%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -10.54%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -10.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.79'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -7.52%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -7.75%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.167'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.86'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -7.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.653.00'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.05%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -7.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.26'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.55%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0920'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -8.78%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -6.54%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '163.72'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -5.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '48.08'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '29.69'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.12%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.75'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -15.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000281'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -8.56%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.11%  '
